$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "article 85 is live": the rolling blog-article window shifts down by one.
# I7 (ser: 82) -> ser: 83
# E7 (ser: 83) -> ser: 84
# C7 (ser: 84) -> ser: 85
# D7 (meetup)  -> unchanged content, only its shared-string slot shuffled upstream

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 83"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 84"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 85"
